$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook for "Galaxy S7 "
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Galaxy S7 "

$newSheet.Range("A1").Value = "Test Parameters"
$newSheet.Range("B1").NumberFormat = "@"
$newSheet.Range("B1").Value = "0.01"

$newSheet.Columns.Item(1).ColumnWidth = 15.60546875
$newSheet.Columns.Item(2).ColumnWidth = 4.859375
